# Rename the "_old"/"_new" suffixed column headers to the respective
# input-file-version suffixes ("_FV2310" / "_FV2404"), freeze the header
# row, and turn the header range into a real Excel Table ("Table1")
# spanning the used range so the autofilter / header formatting travels
# with the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310",
    "diff",
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Freeze the header row (split below row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Promote the A1:U61 range to a real table so the autofilter / header
# band travel with the data.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U61"), $null, 1)
$tbl.Name = "Table1"
